$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Refresh the per-row "time_taken" query timestamps on the "data" sheet ---
# (rows 2..92, column F) to reflect the latest panel re-fetch.
$newQueryTimes = @(
    "2021-10-05 14:21:53.763501",
    "2021-10-05 14:21:53.763510",
    "2021-10-05 14:21:53.763513",
    "2021-10-05 14:21:53.763516",
    "2021-10-05 14:21:53.763519",
    "2021-10-05 14:21:53.763522",
    "2021-10-05 14:21:53.763524",
    "2021-10-05 14:21:53.763527",
    "2021-10-05 14:21:53.763530",
    "2021-10-05 14:21:53.763532",
    "2021-10-05 14:21:53.763535",
    "2021-10-05 14:21:53.763538",
    "2021-10-05 14:21:53.763540",
    "2021-10-05 14:21:53.763543",
    "2021-10-05 14:21:53.763545",
    "2021-10-05 14:21:53.763548",
    "2021-10-05 14:21:53.763551",
    "2021-10-05 14:21:53.763553",
    "2021-10-05 14:21:53.763556",
    "2021-10-05 14:21:53.763559",
    "2021-10-05 14:21:53.763562",
    "2021-10-05 14:21:53.763564",
    "2021-10-05 14:21:53.763567",
    "2021-10-05 14:21:53.763570",
    "2021-10-05 14:21:53.763573",
    "2021-10-05 14:21:53.763576",
    "2021-10-05 14:21:53.763578",
    "2021-10-05 14:21:53.763581",
    "2021-10-05 14:21:53.763584",
    "2021-10-05 14:21:53.763586",
    "2021-10-05 14:21:53.763589",
    "2021-10-05 14:21:53.763592",
    "2021-10-05 14:21:53.763595",
    "2021-10-05 14:21:53.763598",
    "2021-10-05 14:21:53.763600",
    "2021-10-05 14:21:53.763603",
    "2021-10-05 14:21:53.763606",
    "2021-10-05 14:21:53.763608",
    "2021-10-05 14:21:53.763611",
    "2021-10-05 14:21:53.763615",
    "2021-10-05 14:21:53.763618",
    "2021-10-05 14:21:53.763621",
    "2021-10-05 14:21:53.763624",
    "2021-10-05 14:21:53.763627",
    "2021-10-05 14:21:53.763629",
    "2021-10-05 14:21:53.763632",
    "2021-10-05 14:21:53.763635",
    "2021-10-05 14:21:53.763637",
    "2021-10-05 14:21:53.763640",
    "2021-10-05 14:21:53.763642",
    "2021-10-05 14:21:53.763645",
    "2021-10-05 14:21:53.763648",
    "2021-10-05 14:21:53.763651",
    "2021-10-05 14:21:53.763653",
    "2021-10-05 14:21:53.763656",
    "2021-10-05 14:21:53.763659",
    "2021-10-05 14:21:53.763661",
    "2021-10-05 14:21:53.763664",
    "2021-10-05 14:21:53.763667",
    "2021-10-05 14:21:53.763669",
    "2021-10-05 14:21:53.763672",
    "2021-10-05 14:21:53.763675",
    "2021-10-05 14:21:53.763677",
    "2021-10-05 14:21:53.763680",
    "2021-10-05 14:21:53.763684",
    "2021-10-05 14:21:53.763687",
    "2021-10-05 14:21:53.763690",
    "2021-10-05 14:21:53.763693",
    "2021-10-05 14:21:53.763696",
    "2021-10-05 14:21:53.763699",
    "2021-10-05 14:21:53.763702",
    "2021-10-05 14:21:53.763705",
    "2021-10-05 14:21:53.763707",
    "2021-10-05 14:21:53.763710",
    "2021-10-05 14:21:53.763713",
    "2021-10-05 14:21:53.763716",
    "2021-10-05 14:21:53.763720",
    "2021-10-05 14:21:53.763724",
    "2021-10-05 14:21:53.763727",
    "2021-10-05 14:21:53.763730",
    "2021-10-05 14:21:53.763733",
    "2021-10-05 14:21:53.763736",
    "2021-10-05 14:21:53.763738",
    "2021-10-05 14:21:53.763741",
    "2021-10-05 14:21:53.763744",
    "2021-10-05 14:21:53.763747",
    "2021-10-05 14:21:53.763750",
    "2021-10-05 14:21:53.763753",
    "2021-10-05 14:21:53.763755",
    "2021-10-05 14:21:53.763758",
    "2021-10-05 14:21:53.763761"
)
for ($i = 0; $i -lt $newQueryTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Range("F$row").Value = $newQueryTimes[$i]
}

# --- Add the new "metadata" worksheet after the existing "data" sheet. ---
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row (row 1), columns B..G
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Ophthalmological ciliopathies"
$ws.Range("C2").Value = 722

# data_version is stored as text "1.19" (not a number) in the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.19"

$ws.Range("E2").Value = "2021-09-06T14:01:30.870320Z"
$ws.Range("F2").Value = "2021-10-05 14:21:53.760080"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/722/?format=json"

# Formatting: match the bold/bordered/centered header style used on the "data" sheet.
$headerRange = $ws.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$indexCell = $ws.Range("A2")
$indexCell.Font.Bold = $true
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160
$indexCell.Borders.LineStyle = 1

# Keep the "data" sheet as the active tab (matches the original workbook view).
$dataSheet.Activate()
